$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.21973717212677
$ws.Range("B1").Value = 2.568229675292969
$ws.Range("C1").Value = 9.288806915283203
$ws.Range("D1").Value = 2.062599897384644
$ws.Range("E1").Value = 1.186188340187073
